$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price / volume figures. Values are stored as plain text in
# the sheet (Price/Volume columns), so a leading apostrophe forces Excel to
# treat number-looking strings ("211.75") as text instead of coercing them
# to numeric values; resetting the Style back to Normal afterwards avoids
# leaving a stray quote-prefix / text number-format on the cell.
$ws.Range("D2").Value = "'28.492.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.04%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.566.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.21%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.08%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'211.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'0.492"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.38%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.12%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'46.21"
$ws.Range("D8").Style = "Normal"
$ws.Range("E9").Value = "'  +0.03%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -1.89%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0592"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.65%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.0887"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.33%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.790.36"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -2.33%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'1.566.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.37%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  -2.79%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'28.490.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.17%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -3.43%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'62.23"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.91%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'228.10"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.06%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -2.79%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E22").Value = "'  +0.10%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -6.09%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'9.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -3.36%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +6.00%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'150.72"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.29%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -2.22%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -2.79%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -4.10%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +0.04%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -1.88%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'1.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -4.09%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'3.22"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.23%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -3.04%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.392.23"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -2.30%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -0.85%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -3.47%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'2.35"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.84%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +2.39%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -1.07%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -1.58%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E43").Value = "'  +2.03%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.787"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -4.50%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -4.24%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.974"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.39%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'62.87"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -3.38%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'1.703.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.27%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -2.00%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'Cronos"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'0.0525"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.39%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'BabyDogeCoin"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'0.0₆0102"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -4.41%  "
$ws.Range("E51").Style = "Normal"
